$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected data values (revision update from source MV)
$updates = @{
    "E134" = 91.2
    "C138" = 92.95999999999999
    "E139" = 91.41
    "E142" = 95.92
    "E143" = 96.09999999999999
    "E144" = 98.28
    "E145" = 99.87
    "E146" = 94.25
    "E147" = 93.72
    "E148" = 107.2
    "E149" = 97.69
    "E150" = 102.9
    "E151" = 100.31
    "E152" = 100.65
    "E153" = 105.22
    "E154" = 93.8
    "E155" = 107.11
    "E156" = 103.71
    "C157" = 97.13
    "E157" = 105.37
    "E158" = 101.85
    "E159" = 91.48
    "E160" = 99.70999999999999
    "E161" = 100.26
    "C162" = 97.25
    "E162" = 102.64
    "E163" = 99.26000000000001
    "E165" = 105.39
    "E166" = 92.25
    "E167" = 104.17
    "E168" = 104.29
    "E169" = 102.06
    "E170" = 100.59
    "E171" = 93.44
    "C172" = 99.40000000000001
    "E172" = 102.22
    "E173" = 98.43000000000001
    "C174" = 101.95
    "C175" = 101.57
    "E175" = 97.53
    "E176" = 101.19
    "E177" = 98.87
    "E178" = 96.62
    "E179" = 104.62
    "E180" = 102.29
    "C182" = 103.54
    "E183" = 91.2
    "E187" = 98.5
    "C188" = 99.98
    "E189" = 99.84
    "C190" = 100.81
    "E190" = 98.34
    "E191" = 102.32
    "C192" = 98.63
    "E192" = 100.67
    "C193" = 99.43000000000001
    "E193" = 107.29
    "C194" = 99.8
    "E194" = 98.02
    "E195" = 95.15000000000001
    "C196" = 102.52
    "E196" = 103.78
    "E197" = 98.25
    "C198" = 99.41
    "C199" = 99.29000000000001
    "C200" = 98.23999999999999
    "E200" = 98.38
    "C201" = 98.23999999999999
    "E201" = 101.01
    "C202" = 98.68000000000001
    "E202" = 99
    "C203" = 97.09
    "E203" = 96.45999999999999
    "C204" = 99.69
    "E204" = 100.9
    "E205" = 102.19
    "C206" = 98.62
    "E206" = 94.28
    "C207" = 94.36
    "E207" = 89.92
    "C208" = 91.20999999999999
    "C209" = 97.31
    "C210" = 97.68000000000001
    "C211" = 98.02
    "E211" = 94.25
    "C212" = 100.1
    "E212" = 96.53
    "C213" = 101.85
    "C214" = 100.85
    "C215" = 101.02
    "C216" = 101.64
    "E216" = 103.4
    "E217" = 96.63
    "C218" = 101.18
    "E218" = 98.78
    "C219" = 103.94
    "E219" = 91.78
    "C220" = 102.3
    "E220" = 101.54
    "C221" = 101.39
    "E221" = 100.09
    "C222" = 101.92
    "E222" = 101.37
    "C223" = 103.11
    "C224" = 102.01
    "C225" = 101.39
    "E225" = 104.48
    "C226" = 100.51
    "C227" = 100.68
    "C229" = 103.04
    "C230" = 100.49
    "E230" = 100.95
    "C231" = 99.75
    "E231" = 89.16
    "C232" = 100.19
    "C233" = 102.65
    "E233" = 97.45999999999999
    "C234" = 102.81
    "E234" = 101.35
    "C235" = 101.09
    "E235" = 95.06999999999999
    "C236" = 101.63
    "E236" = 101.44
    "C237" = 103.21
    "E237" = 97.73
    "C238" = 101.26
    "C239" = 98.31
    "E239" = 99.89
    "C240" = 101.43
    "E240" = 100.99
    "C241" = 104.87
    "E241" = 100.08
    "C242" = 102.67
    "E242" = 101.81
    "C243" = 101.66
    "E243" = 95.04000000000001
    "C244" = 101.03
    "E244" = 103.51
    "C245" = 98.06999999999999
    "E245" = 94.47
    "C246" = 98.22
    "E246" = 93.58
    "C247" = 97.3
    "E247" = 87.59
    "C248" = 98.01000000000001
    "E248" = 97.39
    "B249" = 99
    "C249" = 98.58
    "E249" = 98.54000000000001
    "C250" = 100.49
    "E250" = 96
    "C251" = 102.16
    "E251" = 103.41
    "C252" = 100.54
    "E252" = 102.74
    "C254" = 101.78
    "E254" = 98.61
    "C255" = 101.59
    "E255" = 93.03
    "C256" = 102.97
    "E256" = 110.87
    "C257" = 103.17
    "E257" = 101.57
    "C258" = 101.4
    "E258" = 102.9
    "C259" = 104.27
    "E259" = 103.16
    "B260" = 102.61
    "C260" = 103.39
    "E260" = 108.41
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# New row 261: 01-08-2021
$ws.Range("A261").NumberFormat = "@"
$ws.Range("A261").Value = "01-08-2021"
$ws.Range("A261").Style = "Normal"
$ws.Range("B261").Value = 103.2
$ws.Range("C261").Value = 101.88
$ws.Range("E261").Value = 110.52
